$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the full data block (rows 3-17, columns A-AY) before making changes.
$srcRange = $ws.Range("A3:AY17")
$vals = $srcRange.Value2

# Mapping of new (1-based) array row -> old (1-based) array row within the snapshot.
# This reproduces the row re-ordering described by the diff (row 3 maps to old row4's
# data, row 4 maps to old row6's data, etc.)
$perm = @{1=2; 2=4; 3=1; 4=15; 5=3; 6=5; 7=6; 8=7; 9=8; 10=9; 11=10; 12=11; 13=12; 14=13; 15=14}

$newVals = New-Object 'object[,]' 15,51
for ($r = 1; $r -le 15; $r++) {
    $oldR = $perm[$r]
    for ($c = 1; $c -le 51; $c++) {
        $newVals[$r-1, $c-1] = $vals[$oldR, $c]
    }
}

# Columns Y,Z,AA,AB (Startdatum/Starttid/Slutdatum/Sluttid) hold date/time values stored
# as plain text. Force text format before writing so Excel doesn't silently convert the
# text back into date serial numbers.
$dateColsRange = $ws.Range("Y3:AB17")
$dateColsRange.NumberFormat = "@"

$ws.Range("A3:AY17").Value2 = $newVals

$dateColsRange.NumberFormat = "General"
